$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 91 and 92 (the en dash "–" cid30373 and em dash "—" cid30374 rows),
# shifting all subsequent rows up by two. The shared strings for "–" and "—"
# are no longer referenced afterwards and drop out of sharedStrings.xml, and the
# last two rows of the table (previously 101/102) disappear since data shifted up.
$ws.Rows("91:92").Delete()

# Leave the selection on the (now-empty, shifted-up) two rows that were just
# removed, matching the state right after performing the deletion, and scroll
# the view back up a bit.
$ws.Application.Goto($ws.Range("A91:XFD92"), $true)
$excel.ActiveWindow.ScrollRow = 70
